$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignUp")

# Replace the old test rows (gopalesh/samara/virpal) with newly added test data
$ws.Range("A2").Value = "Govind Patel"
$ws.Range("B2").Value = "govind@asite.com"

$ws.Range("A3").Value = "Namrata Shah"
$ws.Range("B3").Value = "namrata@asite.com"

$ws.Range("A4").Value = "Gajendra Rathod"
$ws.Range("B4").Value = "gajendra@asite.com"

# Emails in column B pick up an explicit (blank) fill format, matching the
# new cellXfs entry added to the workbook's style table.
$ws.Range("B2:B4").Interior.Pattern = 17

# Last-used cell selection when the workbook was saved
[void]$ws.Range("F12").Select()
